$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date-like text in column F stays text (not auto-converted to a date serial)
$ws.Range("F2:F21").NumberFormat = "@"

# Row 2: Senior AI/ML Developer
$ws.Range("A2").Value = 'Senior AI/ML Developer'
$ws.Range("B2").Value = 'Tata Consultancy Services (TCS)'
$ws.Range("C2").Value = 'Irving, TX, US USA'
$ws.Range("E2").Value = 'Generative AI, RAG, TensorFlow, PyTorch, Azure ML, MLflow, FastAPI, Docker, Kubernetes, CI/CD'
$ws.Range("F2").Value = '2026-02-27'
$ws.Range("G2").Value = 'https://www.indeed.com/viewjob?jk=61c8a3c6f4ba6dbe'
$ws.Range("D2").Value = 18.9

# Row 3: Data Scientist - Growth
$ws.Range("A3").Value = 'Data Scientist - Growth'
$ws.Range("B3").Value = 'Art of Problem Solving Academy'
$ws.Range("C3").Value = 'San Diego, CA, US USA'
$ws.Range("E3").Value = 'Data Scientist, Redshift, BigQuery, Snowflake, BigQuery, Redshift, PySpark, Polars, Dask, Tableau'
$ws.Range("F3").Value = '2026-02-27'
$ws.Range("G3").Value = 'https://www.indeed.com/viewjob?jk=7683f53bd78d48c5'
$ws.Range("D3").Value = 15.6

# Row 4: Gen AI Engineer
$ws.Range("A4").Value = 'Gen AI Engineer'
$ws.Range("B4").Value = 'Tata Consultancy Services (TCS)'
$ws.Range("C4").Value = 'Frederick, MD, US USA'
$ws.Range("E4").Value = 'AI Engineer, Generative AI, LangChain, RAG, Prompt Engineering, TensorFlow, PyTorch, Azure ML, Docker, Kubernetes'
$ws.Range("F4").Value = '2026-02-27'
$ws.Range("G4").Value = 'https://www.indeed.com/viewjob?jk=46d11312afda75ac'
$ws.Range("D4").Value = 14.4

# Row 5: Senior Software Engineer - CICD
$ws.Range("A5").Value = 'Senior Software Engineer - CICD'
$ws.Range("B5").Value = 'Klaviyo'
$ws.Range("C5").Value = 'Boston, MA, US USA'
$ws.Range("E5").Value = 'RAG, Copilot, Kubernetes, CI/CD, Terraform, Git, Kafka, MySQL, Python, SQL'
$ws.Range("F5").Value = '2026-02-27'
$ws.Range("G5").Value = 'https://www.indeed.com/viewjob?jk=f4cd3477fb442eb3'
$ws.Range("D5").Value = 14.4

# Row 6: Machine Learning Software Engineer
$ws.Range("A6").Value = 'Machine Learning Software Engineer'
$ws.Range("B6").Value = 'Qualcomm'
$ws.Range("C6").Value = 'San Diego, CA, US USA'
$ws.Range("E6").Value = 'RAG, TensorFlow, PyTorch, Docker, Kubernetes, CI/CD, Jenkins, Terraform, Git, Python'
$ws.Range("F6").Value = '2026-02-27'
$ws.Range("G6").Value = 'https://www.indeed.com/viewjob?jk=196dc3bc7f2912cf'
$ws.Range("D6").Value = 14.4

# Row 7: Computer Vision Engineer (Contract)
$ws.Range("A7").Value = 'Computer Vision Engineer (Contract)'
$ws.Range("B7").Value = 'Overvak'
$ws.Range("C7").Value = 'US USA'
$ws.Range("E7").Value = 'LangChain, TensorFlow, PyTorch, OpenCV, YOLO, Docker, Kafka, Matplotlib, Python, R'
$ws.Range("F7").Value = '2026-02-27'
$ws.Range("G7").Value = 'https://www.indeed.com/viewjob?jk=25462128ccdcd210'
$ws.Range("D7").Value = 13.3

# Row 8: AI Architect
$ws.Range("A8").Value = 'AI Architect'
$ws.Range("B8").Value = 'Tata Consultancy Services (TCS)'
$ws.Range("C8").Value = 'New York, NY, US USA'
$ws.Range("E8").Value = 'AI Engineer, Data Scientist, LangChain, RAG, TensorFlow, PyTorch, MLflow, Docker, Kubernetes, Python'
$ws.Range("F8").Value = '2026-02-27'
$ws.Range("G8").Value = 'https://www.indeed.com/viewjob?jk=369510cb82d8c982'
$ws.Range("D8").Value = 13.3

# Row 9: Computer Vision Engineer (Contract)
$ws.Range("A9").Value = 'Computer Vision Engineer (Contract)'
$ws.Range("B9").Value = 'Overvak'
$ws.Range("C9").Value = 'US USA'
$ws.Range("E9").Value = 'LangChain, TensorFlow, PyTorch, OpenCV, YOLO, Docker, Kafka, Matplotlib, Python, R'
$ws.Range("F9").Value = '2026-02-27'
$ws.Range("G9").Value = 'https://www.indeed.com/viewjob?jk=f1b4d1a6e61b5656'
$ws.Range("D9").Value = 13.3

# Row 10: Data Engineer
$ws.Range("A10").Value = 'Data Engineer'
$ws.Range("B10").Value = 'Halvik'
$ws.Range("C10").Value = 'Washington, DC, US USA'
$ws.Range("E10").Value = 'Redshift, Synapse, Git, Snowflake, Databricks, Redshift, Kafka, Python, SQL, R'
$ws.Range("F10").Value = '2026-02-27'
$ws.Range("G10").Value = 'https://www.indeed.com/viewjob?jk=89d69e074834d3c5'
$ws.Range("D10").Value = 13.3

# Row 11: AI Enablement Engineer
$ws.Range("A11").Value = 'AI Enablement Engineer'
$ws.Range("B11").Value = 'GALAXY'
$ws.Range("C11").Value = 'Remote, US USA'
$ws.Range("E11").Value = 'Machine Learning Engineer, Generative AI, RAG, LLaMA, Copilot, Pinecone, ChromaDB, Git, Python, R'
$ws.Range("F11").Value = '2026-02-27'
$ws.Range("G11").Value = 'https://www.indeed.com/viewjob?jk=6863ac7bba3e6e79'
$ws.Range("D11").Value = 12.2

# Row 12: Data Science, Intern - Summer 2026, Austin, TX
$ws.Range("A12").Value = 'Data Science, Intern - Summer 2026, Austin, TX'
$ws.Range("B12").Value = 'Visa'
$ws.Range("C12").Value = 'Austin, TX, US USA'
$ws.Range("E12").Value = 'Data Scientist, Generative AI, RAG, Copilot, TensorFlow, PyTorch, Git, Kafka, Python, SQL'
$ws.Range("F12").Value = '2026-02-27'
$ws.Range("G12").Value = 'https://www.indeed.com/viewjob?jk=613a16ffa52d1563'
$ws.Range("D12").Value = 12.2

# Row 13: Senior Analytics Engineer
$ws.Range("A13").Value = 'Senior Analytics Engineer'
$ws.Range("B13").Value = 'Guardian Bikes'
$ws.Range("C13").Value = 'US USA'
$ws.Range("E13").Value = 'RAG, Gemini, Copilot, BigQuery, Git, Snowflake, BigQuery, Python, SQL, R'
$ws.Range("F13").Value = '2026-02-27'
$ws.Range("G13").Value = 'https://www.indeed.com/viewjob?jk=78011a5d6c324570'
$ws.Range("D13").Value = 12.2

# Row 14: Senior Data Engineer/Scientist
$ws.Range("A14").Value = 'Senior Data Engineer/Scientist'
$ws.Range("B14").Value = 'Zendar'
$ws.Range("C14").Value = 'Berkeley, CA, US USA'
$ws.Range("E14").Value = 'RAG, BigQuery, Data Lake, Kubernetes, Terraform, BigQuery, NoSQL, Python, SQL, R'
$ws.Range("F14").Value = '2026-02-01'
$ws.Range("G14").Value = 'https://www.indeed.com/viewjob?jk=0212cfdc5a2493f3'
$ws.Range("D14").Value = 12.2

# Row 15: Autonomous Driving Vehicle Perception Engineer
$ws.Range("A15").Value = 'Autonomous Driving Vehicle Perception Engineer'
$ws.Range("B15").Value = 'Quest Global'
$ws.Range("C15").Value = 'Lansing, MI, US USA'
$ws.Range("E15").Value = 'Data Scientist, RAG, TensorFlow, PyTorch, OpenCV, YOLO, Git, Python, R, Optimization'
$ws.Range("F15").Value = '2026-02-27'
$ws.Range("G15").Value = 'https://www.indeed.com/viewjob?jk=b37b03f34f658304'
$ws.Range("D15").Value = 12.2

# Row 16: Junior Frontend Developer (Creative + Full-Stack)
$ws.Range("A16").Value = 'Junior Frontend Developer (Creative + Full-Stack)'
$ws.Range("B16").Value = 'Sapphire Media LLC'
$ws.Range("C16").Value = 'Scottsdale, AZ, US USA'
$ws.Range("E16").Value = 'RAG, Hugging Face, FastAPI, Docker, Git, MySQL, Python, SQL, R, Java'
$ws.Range("F16").Value = '2026-02-27'
$ws.Range("G16").Value = 'https://www.indeed.com/viewjob?jk=bd450f3d86eff079'
$ws.Range("D16").Value = 11.1

# Row 17: Senior Backend Engineer
$ws.Range("A17").Value = 'Senior Backend Engineer'
$ws.Range("B17").Value = 'Tenjin'
$ws.Range("C17").Value = 'San Francisco, CA, US USA'
$ws.Range("E17").Value = 'RAG, Docker, Kubernetes, CI/CD, GitHub Actions, Git, Kafka, PostgreSQL, SQL, R'
$ws.Range("F17").Value = '2026-02-27'
$ws.Range("G17").Value = 'https://www.indeed.com/viewjob?jk=5ab9bde15f7a07b7'
$ws.Range("D17").Value = 11.1

# Row 18: Senior Platform Engineer
$ws.Range("A18").Value = 'Senior Platform Engineer'
$ws.Range("B18").Value = 'Mambu'
$ws.Range("C18").Value = 'Miami, FL, US USA'
$ws.Range("E18").Value = 'RAG, Kubernetes, CI/CD, Terraform, Git, Power BI, Python, R, Scala, Optimization'
$ws.Range("F18").Value = '2026-02-27'
$ws.Range("G18").Value = 'https://www.indeed.com/viewjob?jk=7d08210831495843'
$ws.Range("D18").Value = 11.1

# Row 19: 2026 Intern, Memory and Personalization (Summer)
$ws.Range("A19").Value = '2026 Intern, Memory and Personalization (Summer)'
$ws.Range("B19").Value = 'Samsung Research America'
$ws.Range("C19").Value = 'Mountain View, CA, US USA'
$ws.Range("E19").Value = 'TensorFlow, PyTorch, Jenkins, GitHub Actions, Git, Python, R, Scala, Optimization'
$ws.Range("F19").Value = '2026-02-27'
$ws.Range("G19").Value = 'https://www.indeed.com/viewjob?jk=038e06c5fa27ee8e'
$ws.Range("D19").Value = 10

# Row 20: Intern - AI Software QA Engineer
$ws.Range("A20").Value = 'Intern - AI Software QA Engineer'
$ws.Range("B20").Value = 'Fremont Bank'
$ws.Range("C20").Value = 'Livermore, CA, US USA'
$ws.Range("E20").Value = 'Generative AI, RAG, Hugging Face, Prompt Engineering, CI/CD, Git, Python, R, Optimization'
$ws.Range("F20").Value = '2026-02-26'
$ws.Range("G20").Value = 'https://www.indeed.com/viewjob?jk=2e20f3e88254de62'
$ws.Range("D20").Value = 10

# Row 21: Software Engineer - Ford Pro Tech
$ws.Range("A21").Value = 'Software Engineer - Ford Pro Tech'
$ws.Range("B21").Value = 'Ford Motor Company'
$ws.Range("C21").Value = 'Spokane Valley, WA, US USA'
$ws.Range("E21").Value = 'RAG, CI/CD, Jenkins, Terraform, NoSQL, Python, SQL, R, Java'
$ws.Range("F21").Value = '2026-02-27'
$ws.Range("G21").Value = 'https://www.indeed.com/viewjob?jk=10daa7f62943c1b3'
$ws.Range("D21").Value = 10

# Remove the two job postings that fell out of today's match set
$ws.Rows("22:23").Delete()
